$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Replace the user "karen.arceo" with "adolfo.cardenas" in column A (rows 2-4)
$ws.Range("A2").Value = "adolfo.cardenas"
$ws.Range("A3").Value = "adolfo.cardenas"
$ws.Range("A4").Value = "adolfo.cardenas"

# Update the selection shown in the worksheet view to E2
$ws.Range("E2").Select()
